# Auto-generated Excel COM-interop script
# Applies the scheduled market-data refresh to the profit-tracking sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 29412076
$ws.Range("I33").Value = 35714644
$ws.Range("K33").Value = 35714644
$ws.Range("M33").Value = -35714415
$ws.Range("H74").Value = 3399.9565
$ws.Range("I74").Value = 3279.9
$ws.Range("J74").Value = 3492.3076
$ws.Range("K74").Value = 3279.9
$ws.Range("L74").Value = 3492.3076
$ws.Range("M74").Value = -2343.9
$ws.Range("N74").Value = -5364.3076
$ws.Range("H76").Value = 6159.0347
$ws.Range("I76").Value = 4330.1177
$ws.Range("J76").Value = 8750
$ws.Range("K76").Value = 4330.1177
$ws.Range("L76").Value = 8750
$ws.Range("M76").Value = -4015.1177
$ws.Range("N76").Value = -9380
$ws.Range("H77").Value = 3399.9565
$ws.Range("I77").Value = 3279.9
$ws.Range("J77").Value = 3492.3076
$ws.Range("K77").Value = 16399.5
$ws.Range("L77").Value = 17461.538
$ws.Range("M77").Value = -11719.5
$ws.Range("N77").Value = -26821.538
$ws.Range("H79").Value = 6159.0347
$ws.Range("I79").Value = 4330.1177
$ws.Range("J79").Value = 8750
$ws.Range("K79").Value = 4330.1177
$ws.Range("L79").Value = 8750
$ws.Range("M79").Value = -3238.1177
$ws.Range("N79").Value = -10934
$ws.Range("H132").Value = 5564.121
$ws.Range("I132").Value = 2779.9443
$ws.Range("J132").Value = 8905.134
$ws.Range("K132").Value = 8339.832900000001
$ws.Range("L132").Value = 26715.402
$ws.Range("M132").Value = -5809.832900000001
$ws.Range("N132").Value = -31775.402
$ws.Range("H135").Value = 23811716
$ws.Range("I135").Value = 1883.5883
$ws.Range("J135").Value = 125003500
$ws.Range("K135").Value = 16952.2947
$ws.Range("L135").Value = 1125031500
$ws.Range("M135").Value = -14417.2947
$ws.Range("N135").Value = -1125036570
$ws.Range("H138").Value = 7144690.5
$ws.Range("I138").Value = 1433.4
$ws.Range("J138").Value = 16669033
$ws.Range("K138").Value = 4300.200000000001
$ws.Range("L138").Value = 50007099
$ws.Range("M138").Value = 839.7999999999993
$ws.Range("N138").Value = -50017379

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 35000
$ws.Range("J23").Value = 35000
$ws.Range("L23").Value = 35000
$ws.Range("N23").Value = -35518
$ws.Range("H27").Value = 3000
$ws.Range("J27").Value = 3000
$ws.Range("L27").Value = 3000
$ws.Range("N27").Value = -3368
$ws.Range("H61").Value = 14709353
$ws.Range("I61").Value = 22730182
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 22730182
$ws.Range("L61").Value = 4500
$ws.Range("M61").Value = -22729970
$ws.Range("N61").Value = -4924
$ws.Range("H63").Value = 250001250
$ws.Range("I63").Value = 500000000
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 500000000
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -499999314
$ws.Range("N63").Value = -3872
$ws.Range("H66").Value = 250001250
$ws.Range("I66").Value = 500000000
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 2500000000
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -2499996568
$ws.Range("N66").Value = -19364
$ws.Range("H74").Value = 11629884
$ws.Range("J74").Value = 3226.8235
$ws.Range("L74").Value = 3226.8235
$ws.Range("N74").Value = -4974.8235
$ws.Range("H77").Value = 11629884
$ws.Range("J77").Value = 3226.8235
$ws.Range("L77").Value = 16134.1175
$ws.Range("N77").Value = -24870.1175
$ws.Range("H88").Value = 2390.923
$ws.Range("I88").Value = 2242.4443
$ws.Range("J88").Value = 2725
$ws.Range("K88").Value = 2242.4443
$ws.Range("L88").Value = 2725
$ws.Range("M88").Value = -1836.4443
$ws.Range("N88").Value = -3537
$ws.Range("H91").Value = 2390.923
$ws.Range("I91").Value = 2242.4443
$ws.Range("J91").Value = 2725
$ws.Range("K91").Value = 2242.4443
$ws.Range("L91").Value = 2725
$ws.Range("M91").Value = -838.4443000000001
$ws.Range("N91").Value = -5533
$ws.Range("H110").Value = 1950
$ws.Range("J110").Value = 2900
$ws.Range("L110").Value = 2900
$ws.Range("N110").Value = -6990
$ws.Range("H132").Value = 9618300
$ws.Range("I132").Value = 11365490
$ws.Range("K132").Value = 34096470
$ws.Range("M132").Value = -34093940
$ws.Range("H136").Value = 14709353
$ws.Range("I136").Value = 22730182
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 68190546
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -68187996
$ws.Range("N136").Value = -18600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 13515126
$ws.Range("I86").Value = 1653.5652
$ws.Range("J86").Value = 35715828
$ws.Range("K86").Value = 1653.5652
$ws.Range("L86").Value = 35715828
$ws.Range("M86").Value = -530.5652
$ws.Range("N86").Value = -35718074
$ws.Range("H89").Value = 13515126
$ws.Range("I89").Value = 1653.5652
$ws.Range("J89").Value = 35715828
$ws.Range("K89").Value = 8267.826000000001
$ws.Range("L89").Value = 178579140
$ws.Range("M89").Value = -2651.826000000001
$ws.Range("N89").Value = -178590372
$ws.Range("H105").Value = 4436.857
$ws.Range("I105").Value = 3165.5557
$ws.Range("J105").Value = 4876.923
$ws.Range("K105").Value = 3165.5557
$ws.Range("L105").Value = 4876.923
$ws.Range("M105").Value = -1418.5557
$ws.Range("N105").Value = -8370.922999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2144.5217
$ws.Range("I58").Value = 1024.0555
$ws.Range("J58").Value = 6178.2
$ws.Range("K58").Value = 1024.0555
$ws.Range("L58").Value = 6178.2
$ws.Range("M58").Value = -821.0554999999999
$ws.Range("N58").Value = -6584.2
$ws.Range("H132").Value = 3218.0417
$ws.Range("I132").Value = 2348.389
$ws.Range("J132").Value = 5827
$ws.Range("K132").Value = 7045.167
$ws.Range("L132").Value = 17481
$ws.Range("M132").Value = -4515.167
$ws.Range("N132").Value = -22541
$ws.Range("H136").Value = 2144.5217
$ws.Range("I136").Value = 1024.0555
$ws.Range("J136").Value = 6178.2
$ws.Range("K136").Value = 3072.1665
$ws.Range("L136").Value = 18534.6
$ws.Range("M136").Value = -522.1664999999998
$ws.Range("N136").Value = -23634.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 305.35715
$ws.Range("I69").Value = 305.35715
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 916.0714499999999
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -105.0714499999999
$ws.Range("H72").Value = 305.35715
$ws.Range("I72").Value = 305.35715
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 2748.21435
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = 1307.78565
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1000510.6
$ws.Range("I3").Value = 1250488.2
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 1250488.2
$ws.Range("L3").Value = 600
$ws.Range("M3").Value = -1250372.2
$ws.Range("N3").Value = -832
$ws.Range("H13").Value = 559
$ws.Range("I13").Value = 399.66666
$ws.Range("J13").Value = 798
$ws.Range("K13").Value = 399.66666
$ws.Range("L13").Value = 798
$ws.Range("M13").Value = -260.66666
$ws.Range("N13").Value = -1076
$ws.Range("H14").Value = 5000068
$ws.Range("I14").Value = 5000068
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 5000068
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -4999900
$ws.Range("H70").Value = 17786.227
$ws.Range("I70").Value = 27947.44
$ws.Range("J70").Value = 4416.2104
$ws.Range("K70").Value = 27947.44
$ws.Range("L70").Value = 4416.2104
$ws.Range("M70").Value = -27677.44
$ws.Range("N70").Value = -4956.2104
$ws.Range("H73").Value = 17786.227
$ws.Range("I73").Value = 27947.44
$ws.Range("J73").Value = 4416.2104
$ws.Range("K73").Value = 27947.44
$ws.Range("L73").Value = 4416.2104
$ws.Range("M73").Value = -27011.44
$ws.Range("N73").Value = -6288.2104
$ws.Range("H80").Value = 19426174
$ws.Range("I80").Value = 25643596
$ws.Range("J80").Value = 5955096.5
$ws.Range("K80").Value = 25643596
$ws.Range("L80").Value = 5955096.5
$ws.Range("M80").Value = -25642598
$ws.Range("N80").Value = -5957092.5
$ws.Range("H83").Value = 19426174
$ws.Range("I83").Value = 25643596
$ws.Range("J83").Value = 5955096.5
$ws.Range("K83").Value = 128217980
$ws.Range("L83").Value = 29775482.5
$ws.Range("M83").Value = -128212988
$ws.Range("N83").Value = -29785466.5
$ws.Range("H132").Value = 5078.7
$ws.Range("I132").Value = 3569.6667
$ws.Range("J132").Value = 6084.722
$ws.Range("K132").Value = 10709.0001
$ws.Range("L132").Value = 18254.166
$ws.Range("M132").Value = -8179.000100000001
$ws.Range("N132").Value = -23314.166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 28776
$ws.Range("J29").Value = 28776
$ws.Range("L29").Value = 28776
$ws.Range("N29").Value = -29366
$ws.Range("H132").Value = 8339704
$ws.Range("J132").Value = 29423902
$ws.Range("L132").Value = 88271706
$ws.Range("N132").Value = -88276766
$ws.Range("H134").Value = 65000
$ws.Range("J134").Value = 65000
$ws.Range("L134").Value = 65000
$ws.Range("N134").Value = -75140
$ws.Range("H135").Value = 113700
$ws.Range("J135").Value = 113700
$ws.Range("L135").Value = 113700
$ws.Range("N135").Value = -123840
$ws.Range("H136").Value = 35729764
$ws.Range("I136").Value = 55560856
$ws.Range("J136").Value = 33802
$ws.Range("K136").Value = 166682568
$ws.Range("L136").Value = 101406
$ws.Range("M136").Value = -166680018
$ws.Range("N136").Value = -106506

Write-Output "Sheets updated via scheduled runner"
